# Fruta / hortaliza, semanal
# Insert a new weekly record at row 108 (pushing the existing rows 108-143
# down to 109-144) on the active sheet of the "Poroto verde" workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 108; everything currently at/after
# row 108 shifts down by one (new dimension becomes A1:R144).
$ws.Rows("108").Insert()

# Populate the newly inserted row 108 with the new weekly observation.
$ws.Cells.Item(108, 1).Value = 7
$ws.Cells.Item(108, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(108, 3).Value = "Ñuble"
$ws.Cells.Item(108, 4).Value = 45093
$ws.Cells.Item(108, 5).Value = 16
$ws.Cells.Item(108, 6).Value = 100112031
$ws.Cells.Item(108, 7).Value = "Poroto verde"
$ws.Cells.Item(108, 8).Value = "Magnum"
$ws.Cells.Item(108, 9).Value = "Primera"
$ws.Cells.Item(108, 10).Value = 30
$ws.Cells.Item(108, 11).Value = 25000
$ws.Cells.Item(108, 12).Value = 25000
$ws.Cells.Item(108, 13).Value = 25000
$ws.Cells.Item(108, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(108, 15).Value = "Perú"
$ws.Cells.Item(108, 16).Value = 1000
$ws.Cells.Item(108, 17).Value = 25
$ws.Cells.Item(108, 18).Value = "Hortaliza"
